$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 68301
$ws.Range("B2").Value = "Davi Rodrigues"
$ws.Range("C2").Value = "Marketing"
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 45104
$ws.Range("G2").Value = 10674.45

# Row 3
$ws.Range("A3").Value = 10839
$ws.Range("B3").Value = "Dr. Pedro Ferreira"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45083
$ws.Range("G3").Value = 3686.3

# Row 4
$ws.Range("A4").Value = 71524
$ws.Range("B4").Value = "Clara da Mata"
$ws.Range("C4").Value = "Operações"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 3600.45

# Row 5
$ws.Range("A5").Value = 72454
$ws.Range("B5").Value = "João Lucas Azevedo"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45088
$ws.Range("G5").Value = 10480.21

# Row 6
$ws.Range("A6").Value = 62267
$ws.Range("B6").Value = "Lara Rodrigues"
$ws.Range("C6").Value = "Marketing"
$ws.Range("D6").Value = "Doença"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 3078.16

# Row 7
$ws.Range("A7").Value = 43152
$ws.Range("B7").Value = "Murilo Pereira"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45080
$ws.Range("G7").Value = 4106.91

# Row 8
$ws.Range("A8").Value = 14644
$ws.Range("B8").Value = "Srta. Laís da Mata"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45078
$ws.Range("G8").Value = 6839.2

# Row 9
$ws.Range("A9").Value = 92199
$ws.Range("B9").Value = "Stephany Nogueira"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("F9").Value = 45099
$ws.Range("G9").Value = 10175.57

# Row 10
$ws.Range("A10").Value = 53316
$ws.Range("B10").Value = "Luiz Gustavo Peixoto"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Doença"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45103
$ws.Range("G10").Value = 8695.440000000001

# Row 11
$ws.Range("A11").Value = 90960
$ws.Range("B11").Value = "Kamilly Martins"
$ws.Range("C11").Value = "TI"
$ws.Range("D11").Value = "Doença"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 4372.14
